## Remove messageBox useless / adjustement file name
## -> update the two "0€" placeholder amounts (Repas / Hébergement)
##    in the first data row of the expense table to the new sample
##    values used for testing (999€ and 222€).

$d = $word.ActiveDocument

$tbl = $d.Tables.Item(1)

# Row 2 = first data row ("01/01/2018").
# Column 7 = "Repas"       -> 0€  becomes 999€
# Column 8 = "Hébergement" -> 0€  becomes 222€

$repasCell = $tbl.Cell(2, 7)
$repasRange = $repasCell.Range
$repasRange.End = $repasRange.End - 1
$repasRange.Text = "999€"

$hebergementCell = $tbl.Cell(2, 8)
$hebergementRange = $hebergementCell.Range
$hebergementRange.End = $hebergementRange.End - 1
$hebergementRange.Text = "222€"
